# "fixing unit id so we can demo the hospital analysis page"
#
# The Dictionary sheet's raw->standard column mapping had row 16 pointing
# "Pruvodka.Oddeleni.Zkratka" at "HospitalId". Split that into two rows:
#   row 16: Pruvodka.Oddeleni.Zkratka -> UnitId      (fixed mapping)
#   row 17 (new): hosp_id             -> HospitalId  (new mapping, inserted)
# Every row from the old row 17 onward shifts down by one.

$wb = $excel.ActiveWorkbook

$dict = $wb.Worksheets.Item("Dictionary")

# Insert a new row at 17, pushing the old row 17 (and everything below) down.
$dict.Rows.Item(17).Insert()

# Row 16 keeps its raw_column_name, but now maps to the new UnitId standard name.
$dict.Range("B16").Value = "UnitId"

# The freshly inserted row 17 carries the old HospitalId mapping under the
# new raw column name "hosp_id".
$dict.Range("A17").Value = "hosp_id"
$dict.Range("B17").Value = "HospitalId"

# Restore view state: Config sheet scrolled/selected at C32, then leave focus
# back on Dictionary at A18 (matches the saved selection state in the diff).
$cfg = $wb.Worksheets.Item("Config")
$cfg.Activate()
$cfg.Range("C32").Select()

$dict.Activate()
$dict.Range("A18").Select()
